$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update scattered "B" column (sheet column C) imputed/missing values ---
$ws.Range("C6").Value = 15.1
$ws.Range("C8").ClearContents()
$ws.Range("C18").Value = 11.5
$ws.Range("C20").ClearContents()
$ws.Range("C23").Value = 12.2
$ws.Range("C25").ClearContents()

# --- Remove the "RM 232" row (row 26) entirely ---
$ws.Rows(26).Delete()

# After the delete above, the former "SC 92" row has shifted up to row 27.
# Remove it entirely as well.
$ws.Rows(27).Delete()

# --- Fix up the "A" column (sheet column B) values that moved around ---
$ws.Range("B27").Value = -20.4
$ws.Range("B28").ClearContents()
$ws.Range("B29").ClearContents()
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("B32").ClearContents()
